$d = $word.ActiveDocument

# --- Change 1: "Written: Spring Term 2020" -> "Written: Autumn Term 2020" ---
# The source run holds "Spring Term 2020" as one run; the target keeps that run's
# formatting on a new "Autumn" run and puts the remainder (" Term 2020") into a
# second, separately-formatted run. Re-stamping the (unchanged) font values on the
# replaced range forces Word to keep it as its own run instead of re-merging with
# the neighbouring, identically formatted " Term 2020" text.
$rng = $d.Content
$rng.Find.Execute("Spring") | Out-Null
$rng.Text = "Autumn"
$rng.Font.Name = "Candara"
$rng.Font.NameAscii = "Candara"
$rng.Font.NameBi = "Arial"
$rng.Font.NameFarEast = "Candara"

# --- Change 2: "Date: 9th April 2020" -> "Date: 31st August 2020" ---
$dateRng = $d.Content
$dateRng.Find.Execute("Date: 9th April 2020") | Out-Null
$base = $dateRng.Start

# "9" -> "31"
$d.Range($base + 6, $base + 7).Text = "31"

# "th" -> "st" (offset shifted by +1 because "31" is one character longer than "9")
$d.Range($base + 8, $base + 10).Text = "st"

# " April" -> " August" (the space run and the "April" run collapse into one run)
$d.Range($base + 10, $base + 16).Text = " August"
